$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 779, shifting all existing rows (779-820) down to (780-821).
$ws.Rows.Item(779).Insert()

# Populate the newly-inserted row 779 with the new data point. Column A holds a
# date-formatted string ("YYYY/MM/DD") stored as plain text elsewhere in the
# sheet, so force text formatting before assigning it to avoid Excel silently
# reinterpreting it as a date serial number, then clear the formatting back so
# the cell carries no explicit style (matching the rest of the data rows).
$ws.Cells.Item(779, 1).NumberFormat = "@"
$ws.Cells.Item(779, 1).Value = "2026/02/06"
$ws.Cells.Item(779, 1).ClearFormats()

$ws.Cells.Item(779, 2).Value = "金"
$ws.Cells.Item(779, 3).Value = 11
$ws.Cells.Item(779, 4).Value = 201
